$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZK_Bug_Report")
$ws.Activate()

# --- Content change: populate the "Date Fixed" field for each of the
# four bug-report blocks on the sheet (rows 16, 34, 52, 70) with the
# completion date 2023-12-04 (serial 45264), matching the date format
# already used by the "Date Reported" cells (e.g. B3). ---
$xlPasteFormats = -4122

$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("B16").Value = 45264

$ws.Range("B21").Copy()
$ws.Range("B34").PasteSpecial($xlPasteFormats)
$ws.Range("B34").Value = 45264

$ws.Range("B39").Copy()
$ws.Range("B52").PasteSpecial($xlPasteFormats)
$ws.Range("B52").Value = 45264

$ws.Range("B57").Copy()
$ws.Range("B70").PasteSpecial($xlPasteFormats)
$ws.Range("B70").Value = 45264

$excel.CutCopyMode = 0

# --- View state: zoom in on the sheet and move the selection/scroll
# position down to the last bug-report block. ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("B82").Select()
